# biometric_type.xlsx -> add French ("fra") translations for the three
# existing biometric types (Fingerprint/FNR, Iris/IRS, Photo/PHT), mirroring
# the English rows already present, and drop the trailing blank rows left
# over at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the 4 leftover blank rows at the very bottom of the sheet
# (1048573:1048576) before we start inserting, so row numbers for the
# inserts below stay simple.
$ws.Range("A1048573:A1048576").EntireRow.Delete()

# Row 5: fra / FNR / Empreintes digitales / Empreintes digitales du demandeur / TRUE
# (copy row 2 - eng/FNR/Fingerprint/.../TRUE - so style/number formats match)
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = "fra"
$ws.Range("C5").Value = "Empreintes digitales"
$ws.Range("D5").Value = "Empreintes digitales du demandeur"

# Row 6: fra / IRS / Iris / Iris du demandeur / TRUE
# (copy row 3 - eng/IRS/Iris/.../TRUE)
$ws.Rows.Item(3).Copy()
$ws.Rows.Item(6).Insert()
$ws.Range("A6").Value = "fra"
$ws.Range("D6").Value = "Iris du demandeur"

# Row 7: fra / PHT / Photo / Photo du visage du demandeur / TRUE
# (copy row 4 - eng/PHT/Photo/.../TRUE)
$ws.Rows.Item(4).Copy()
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = "fra"
$ws.Range("D7").Value = "Photo du visage du demandeur"

# Match the saved selection / cursor position from the authored workbook.
$ws.Range("G6").Select() | Out-Null

# Cosmetic default-column-width tweak carried by the original commit
# (best-effort; harmless if the host doesn't persist it).
$ws.StandardWidth = 8.4296875
